# Applies the price/volume refresh described in the commit:
# "Updated cryptos list on Sun Sep 24 10:45:04 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $cell = $ws.Range($cellRef)
    # The source data is plain text (coinranking.com scrape formatted as
    # "12.345.67"-style thousand-grouped strings / "  +0.12%  " strings).
    # Pre-format as Text only when the new value would otherwise be parsed
    # as a real number (e.g. "1.01", "65.08") so it stays a literal string
    # like the rest of the column, instead of being auto-converted.
    if ($text -match "^[+-]?\d+(\.\d+)?$") {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

Set-TextValue "D2" "26.730.04"
Set-TextValue "E2" "  +0.24%  "
Set-TextValue "D3" "1.601.86"
Set-TextValue "E3" "  +0.25%  "
Set-TextValue "E4" "  +0.11%  "
Set-TextValue "D5" "211.83"
Set-TextValue "E5" "  +0.15%  "
Set-TextValue "E6" "  -0.15%  "
Set-TextValue "D7" "1.01"
Set-TextValue "E7" "  +0.17%  "
Set-TextValue "E8" "  +0.12%  "
Set-TextValue "E9" "  -0.08%  "
Set-TextValue "E10" "  +0.67%  "
Set-TextValue "E11" "  +0.80%  "
Set-TextValue "D12" "1.825.98"
Set-TextValue "E12" "  +0.21%  "
Set-TextValue "D13" "1.606.39"
Set-TextValue "E13" "  +0.48%  "
Set-TextValue "E14" "  +0.57%  "
Set-TextValue "D15" "0.524"
Set-TextValue "E15" "  -0.01%  "
Set-TextValue "D16" "65.08"
Set-TextValue "E16" "  -0.15%  "
Set-TextValue "D17" "0.0₃0737"
Set-TextValue "E17" "  +0.17%  "
Set-TextValue "D18" "210.00"
Set-TextValue "E18" "  +0.13%  "
Set-TextValue "D19" "1.01"
Set-TextValue "E19" "  +0.15%  "
Set-TextValue "E20" "  +1.83%  "
Set-TextValue "D21" "4.29"
Set-TextValue "E21" "  -0.18%  "
Set-TextValue "E22" "  -3.08%  "
Set-TextValue "D23" "8.99"
Set-TextValue "E23" "  +0.11%  "
Set-TextValue "D24" "143.62"
Set-TextValue "E24" "  -0.46%  "
Set-TextValue "E25" "  +0.16%  "
Set-TextValue "E27" "  -0.83%  "
Set-TextValue "E28" "  +0.36%  "
Set-TextValue "E29" "  -1.16%  "
Set-TextValue "D30" "1.16"
Set-TextValue "E30" "  +0.48%  "
Set-TextValue "E31" "  +0.34%  "
Set-TextValue "E32" "  -0.09%  "
Set-TextValue "D33" "1.289.97"
Set-TextValue "E33" "  +0.13%  "
Set-TextValue "E34" "  +0.81%  "
Set-TextValue "D35" "1.50"
Set-TextValue "E35" "  +0.31%  "
Set-TextValue "D36" "0.602"
Set-TextValue "E36" "  -2.81%  "
Set-TextValue "E37" "  +10.28%  "
Set-TextValue "E38" "  -0.05%  "
Set-TextValue "E39" "  -0.48%  "
Set-TextValue "D40" "5.40"
Set-TextValue "E40" "  -2.03%  "
Set-TextValue "E41" "  -0.32%  "
Set-TextValue "D42" "0.785"
Set-TextValue "E42" "  -0.01%  "
Set-TextValue "D43" "62.87"
Set-TextValue "E43" "  -1.02%  "
Set-TextValue "D44" "1.738.39"
Set-TextValue "E44" "  +0.16%  "
Set-TextValue "D45" "90.55"
Set-TextValue "E45" "  -0.26%  "
Set-TextValue "E46" "  -1.39%  "
Set-TextValue "E47" "  +0.22%  "
Set-TextValue "E48" "  +1.50%  "
Set-TextValue "E49" "  +0.23%  "
Set-TextValue "D50" "7.43"
Set-TextValue "E50" "  +0.17%  "
Set-TextValue "E51" "  +0.89%  "
